$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Row 4 - fill in day 3 (F4) and day 4 (G4) hours
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 4

# Row 5 - fill in day 2 (E5), day 3 (F5) and day 4 (G5) hours
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 2

# Row 6 - new task "oef database les 3" with initial estimate 5 and day 4 (G6) hours
$ws.Range("C6").Value = "oef database les 3"
$ws.Range("D6").Value = 5
$ws.Range("G6").Value = 6

# Move the active selection to E7 (matches the author's final cursor position)
$ws.Range("E7").Select()
